$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New formulas added below row 17: F19 = B7-C7, G19 = F19/C7 (percentage)
$ws.Range("F19").Formula = "=B7-C7"
$ws.Range("G19").Formula = "=F19/C7"
$ws.Range("G19").NumberFormat = "0.00%"

$ws.Range("G19").Select()
